# Fix albedo table overlap on slide 10 in G7_C3_W2 PPTX
# Reduce table row heights from 304800 EMU (24pt) to 200000 EMU (~15.75pt)
# so the albedo reference table no longer overlaps the WORD BANK box below it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)

# Locate the "Table 20" shape (the albedo reference data table) on the slide.
$tableShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tableShape = $candidate
    }
}

$tbl = $tableShape.Table

# PowerPoint's COM object model expresses row heights in points;
# 200000 EMU == 200000 / 12700 points.
$newRowHeightPts = 200000 / 12700

for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    $tbl.Rows.Item($r).Height = $newRowHeightPts
}
